# ERP-514 - Remove , from the court address
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove trailing commas from the tribunal address values (Manchester & Glasgow)
$used = $ws.UsedRange
$used.Replace(",", "")

# Move selection to B5 (as recorded in the saved view state)
$ws.Range("B5").Select()
